$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.620.38'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '2.386.51'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = "'503.30"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').Value = "'131.89"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D9').Value = '2.390.86'
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').Value = "'0.0972"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').Value = "'0.150"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = "'0.322"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = "'4.64"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').Value = '2.813.96'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '56.554.84'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = "'21.61"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.25%  '
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').Value = '2.360.23'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').Value = "'308.32"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').Value = "'6.27"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.27%  '
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').Value = "'5.55"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.97%  '
$ws.Range('D25').Value = "'67.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').Value = "'0.376"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.88%  '
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').Value = "'175.31"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').Value = '0.0₃0722'
$ws.Range('E31').Value = '  +0.88%  '
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('D34').Value = "'5.85"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.82%  '
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').Value = "'0.999"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = "'17.83"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('D39').Value = "'3.80"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.73%  '
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D40').Value = "'0.824"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.84%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = "'36.79"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.60%  '
$ws.Range('E42').Value = '  -0.26%  '
$ws.Range('D43').Value = "'131.15"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.86%  '
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = "'4.84"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').Value = "'250.48"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').Value = "'0.0907"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = "'0.0483"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.19%  '
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').Value = "'17.02"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.32%  '
